$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column G entirely ("CONSTRUCTION AND OTHERS"). This removes the
# now-unused shared string, shifts column H ("NON-ENERGY CONSUMPTION") and
# I ("TOT") one column to the left, and shrinks the used range to A1:H74.
$ws.Columns("G").Delete()

# The former TOT column (now column H) was originally computed including the
# "CONSTRUCTION AND OTHERS" figures; after removing that category the totals
# were recalculated to exclude it. Write the corrected TOT values back in.
$totals = @{
    2 = 58.83; 3 = 61.69; 4 = 66; 5 = 71.69; 6 = 76.19; 7 = 78.47;
    8 = 83.20999999999999; 9 = 85; 10 = 88.91; 11 = 93.78; 12 = 94.53;
    13 = 91.23; 14 = 92.33; 15 = 92.08; 16 = 95.22; 17 = 99.63; 18 = 105.25;
    19 = 108.7; 20 = 109.66; 21 = 111.4; 22 = 114.04; 23 = 115.48; 24 = 116.71;
    25 = 120; 26 = 126.26; 27 = 131.52; 28 = 138.21; 29 = 142.71; 30 = 147.94;
    31 = 155.22; 32 = 156.81; 33 = 156.52; 34 = 159.68; 35 = 161.36; 36 = 170.3;
    37 = 173.11; 38 = 178.28; 39 = 189.14; 40 = 195.41; 41 = 190.81; 42 = 207.71;
    43 = 214.86; 44 = 229.77; 45 = 233.71; 46 = 237.76; 47 = 232.54; 48 = 228.65;
    49 = 233; 50 = 228.67; 51 = 230.3; 52 = 226.96; 53 = 237.19; 54 = 237.73;
    55 = 243.2; 56 = 244.02; 57 = 249.46; 58 = 249.85; 59 = 254.33; 60 = 257.85;
    61 = 259.93; 62 = 263.49; 63 = 266.5; 64 = 269.14; 65 = 272.88; 66 = 275.36;
    67 = 277.98; 68 = 281.04; 69 = 284.41; 70 = 286.47; 71 = 289.99; 72 = 292.13;
    73 = 295.52; 74 = 298.58
}

foreach ($row in $totals.Keys) {
    $ws.Range("H$row").Value = $totals[$row]
}
